$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.371.19"
$ws.Range("E2").Value = "  +6.54%  "

$ws.Range("D3").Value = "1.816.30"
$ws.Range("E3").Value = "  +6.48%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "'343.10"
$ws.Range("E5").Value = "  +3.99%  "

$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").Value = "'0.3864"
$ws.Range("E7").Value = "  +5.18%  "

$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'50.49"
$ws.Range("E8").Value = "  +4.34%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3550"
$ws.Range("E9").Value = "  +7.73%  "

$ws.Range("D10").Value = "'1.249"
$ws.Range("E10").Value = "  +7.22%  "

$ws.Range("D11").Value = "'0.07819"
$ws.Range("E11").Value = "  +6.82%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'22.72"
$ws.Range("E12").Value = "  +14.18%  "

$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "'1.002"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").Value = "'6.681"
$ws.Range("E14").Value = "  +7.93%  "

$ws.Range("D15").Value = "'7.281"
$ws.Range("E15").Value = "  +6.96%  "

$ws.Range("D16").Value = "1.814.51"
$ws.Range("E16").Value = "  +6.61%  "

$ws.Range("D17").Value = "'0.00001137"
$ws.Range("E17").Value = "  +6.25%  "

$ws.Range("D18").Value = "'0.06779"
$ws.Range("E18").Value = "  +2.76%  "

$ws.Range("D19").Value = "'87.16"
$ws.Range("E19").Value = "  +7.88%  "

$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.53%  "

$ws.Range("D21").Value = "'18.07"
$ws.Range("E21").Value = "  +12.18%  "

$ws.Range("D22").Value = "'6.630"
$ws.Range("E22").Value = "  +9.92%  "

$ws.Range("D23").Value = "'13.19"
$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("D24").Value = "27.364.99"
$ws.Range("E24").Value = "  +6.62%  "

$ws.Range("D25").Value = "'2.475"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D26").Value = "'2.737"
$ws.Range("E26").Value = "  +10.57%  "

$ws.Range("D27").Value = "'22.17"
$ws.Range("E27").Value = "  +15.82%  "

$ws.Range("D28").Value = "'1.521"
$ws.Range("E28").Value = "  +19.39%  "

$ws.Range("D29").Value = "'154.35"
$ws.Range("E29").Value = "  +3.32%  "

$ws.Range("D30").Value = "2.016.11"
$ws.Range("E30").Value = "  +6.57%  "

$ws.Range("D31").Value = "'137.84"
$ws.Range("E31").Value = "  +7.77%  "

$ws.Range("D32").Value = "'6.494"
$ws.Range("E32").Value = "  +9.19%  "

$ws.Range("D33").Value = "'4.131"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("D34").Value = "'13.91"
$ws.Range("E34").Value = "  +9.82%  "

$ws.Range("D35").Value = "'0.08870"
$ws.Range("E35").Value = "  +4.75%  "

$ws.Range("D36").Value = "'1.726"
$ws.Range("E36").Value = "  +2.26%  "

$ws.Range("D37").Value = "'5.705"
$ws.Range("E37").Value = "  +7.80%  "

$ws.Range("D38").Value = "'0.7025"
$ws.Range("E38").Value = "  +15.42%  "

$ws.Range("D39").Value = "'0.06601"
$ws.Range("E39").Value = "  +6.19%  "

$ws.Range("D40").Value = "'0.02439"
$ws.Range("E40").Value = "  +8.13%  "

$ws.Range("D41").Value = "'0.2283"
$ws.Range("E41").Value = "  +7.87%  "

$ws.Range("D42").Value = "'9.075"
$ws.Range("E42").Value = "  +7.37%  "

$ws.Range("D43").Value = "'1.266"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").Value = "'15.12"
$ws.Range("E44").Value = "  +7.80%  "

$ws.Range("D45").Value = "'0.6647"
$ws.Range("E45").Value = "  +13.76%  "

$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D47").Value = "'3.965"
$ws.Range("E47").Value = "  +3.45%  "

$ws.Range("D48").Value = "'2.197"
$ws.Range("E48").Value = "  +10.02%  "

$ws.Range("D49").Value = "'133.77"

$ws.Range("D50").Value = "'0.07348"
$ws.Range("E50").Value = "  +1.76%  "

$ws.Range("D51").Value = "'81.37"
$ws.Range("E51").Value = "  +6.80%  "
